$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# E2 previously held "Revista Via Iuris" — replace with the index.comunicación link
$ws.Range("E2").Value = "\href{https://indexcomunicacion.es/}{index.comunicación}"

# E3 previously held plain "Estudios sobre el Mensaje Periodístico" — wrap it with its journal link
$ws.Range("E3").Value = "\href{https://revistas.ucm.es/index.php/esmp/index}{Estudios sobre el Mensaje Periodístico}"

# New row 4: add the Via Iuris journal link
$ws.Range("E4").Value = "\href{https://revistas.libertadores.edu.co/index.php/ViaIuris}{Via Iuris}"

# Update the selection to reflect the new active cell
$ws.Range("E4").Select()
